$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove cells D, AN, AP for rows 2 and 3 (AN/AP contents are dropped; AO gets new value)
$ws.Range("D2").ClearContents()
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# Update changed numeric values for rows 2 and 3
$ws.Range("I2").Value = 1.38074398249453
$ws.Range("J2").Value = 1.38074398249453
$ws.Range("K2").Value = -2.41
$ws.Range("L2").Value = -5.273522975929978
$ws.Range("U2").Value = 0.011
$ws.Range("V2").Value = 0.006010928961748633
$ws.Range("W2").Value = -0.2231481481481482
$ws.Range("X2").Value = 0.04158987782926457
$ws.Range("Y2").Value = -0.2647380259774127
$ws.Range("Z2").Value = 0.04330522126409551
$ws.Range("AA2").Value = 0.05979342367099402
$ws.Range("AB2").Value = 0.03189593031297479
$ws.Range("AC2").Value = 0.02789749335801923
$ws.Range("AD2").Value = 6.04
$ws.Range("AF2").Value = 6.04
$ws.Range("AG2").Value = 6.029
$ws.Range("AH2").Value = 0.7674714104193139
$ws.Range("AI2").Value = 0.4649730561970747
$ws.Range("AJ2").Value = 0.7671459473215422
$ws.Range("AK2").Value = 0.4645196085985053
$ws.Range("AL2").Value = 0.702
$ws.Range("AM2").Value = 0.2379999999999999
$ws.Range("AO2").Value = 0.8988603988603989
$ws.Range("AQ2").Value = 2.651260504201681

$ws.Range("I3").Value = 1.38074398249453
$ws.Range("J3").Value = 1.38074398249453
$ws.Range("K3").Value = -2.41
$ws.Range("L3").Value = -5.273522975929978
$ws.Range("U3").Value = 0.011
$ws.Range("V3").Value = 0.006010928961748633
$ws.Range("W3").Value = -0.2231481481481482
$ws.Range("X3").Value = 0.04158987782926457
$ws.Range("Y3").Value = -0.2647380259774127
$ws.Range("Z3").Value = 0.04330522126409551
$ws.Range("AA3").Value = 0.05979342367099402
$ws.Range("AB3").Value = 0.03189593031297479
$ws.Range("AC3").Value = 0.02789749335801923
$ws.Range("AD3").Value = 6.04
$ws.Range("AF3").Value = 6.04
$ws.Range("AG3").Value = 6.029
$ws.Range("AH3").Value = 0.7674714104193139
$ws.Range("AI3").Value = 0.4649730561970747
$ws.Range("AJ3").Value = 0.7671459473215422
$ws.Range("AK3").Value = 0.4645196085985053
$ws.Range("AL3").Value = 0.702
$ws.Range("AM3").Value = 0.2379999999999999
$ws.Range("AO3").Value = 0.8988603988603989
$ws.Range("AQ3").Value = 2.651260504201681

